# Actualización automática 2025-06-10 17:00:08
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 2 (OTROS): update VENTA and POR CUMPLIR
$ws.Range("D2").Value = 1711.4
$ws.Range("E2").Value = -1711.4

# Row 4 (TOTAL): recompute VENTA, POR CUMPLIR and CUMPLIMIENTO
$ws.Range("D4").Value = 1735.71
$ws.Range("E4").Value = 15764.29
$ws.Range("F4").Value = 0.09918342857142858
